$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column to be treated as text so that
# numeric-looking strings (e.g. "214.50") keep their exact
# formatting instead of being coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.784.40'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.645.26'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").Value = '216.49'
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").Value = '19.22'
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '1.633.67'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D14").Value = '0.529'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '64.71'
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").Value = '26.769.95'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '214.50'
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").Value = '4.39'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("E21").Value = '  +12.45%  '
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E23").Value = '  -1.62%  '
$ws.Range("D24").Value = '146.90'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").Value = '7.18'
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").Value = '3.36'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").Value = '3.01'
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("D33").Value = '1.291.13'
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").Value = '1.54'
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("D36").Value = '0.0176'
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("D38").Value = '0.823'
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("D40").Value = '0.806'
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("D42").Value = '5.32'
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("D43").Value = '1.783.82'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").Value = '61.95'
$ws.Range("E44").Value = '  +3.45%  '
$ws.Range("D45").Value = '91.74'
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("E46").Value = '  +0.91%  '
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("D49").Value = '7.64'
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").Value = '0.0973'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '0.407'
$ws.Range("E51").Value = '  +0.10%  '

# Restore the default cell style on the Price column so the
# workbook's styling matches the original (un-styled) cells.
$ws.Range("D2:D51").Style = "Normal"
